$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.024.41"
$ws.Range("E2").Value = "  -6.12%  "
$ws.Range("D3").Value = "2.885.39"
$ws.Range("E3").Value = "  -3.57%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.36%  "
$ws.Range("D8").Value = "2.886.17"
$ws.Range("E8").Value = "  -3.40%  "
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000209"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "3.365.98"
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("D17").Value = "2.887.07"
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.20%  "
$ws.Range("D19").Value = "57.118.07"
$ws.Range("E19").Value = "  -6.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "402.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.671"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0982"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.909"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("E36").Value = "  -12.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("E38").Value = "  +6.73%  "
$ws.Range("D39").Value = "0.0₃0619"
$ws.Range("E39").Value = "  -6.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0338"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.54%  "
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").Value = "2.624.33"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "357.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.41%  "
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "119.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.228"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.76%  "
